$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '72.198.39'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.617.48'
$ws.Range('E3').Value = '  +4.05%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '603.63'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '178.06'
$ws.Range('E6').Value = '  +0.96%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  +9.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '2.615.26'
$ws.Range('E10').Value = '  +3.97%  '
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.106.79'
$ws.Range('E14').Value = '  +4.49%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000186'
$ws.Range('E15').Value = '  +4.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '72.070.65'
$ws.Range('E16').Value = '  +3.84%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.48'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.621.27'
$ws.Range('E18').Value = '  +5.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '381.45'
$ws.Range('E19').Value = '  +5.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.55'
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.91'
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.01'
$ws.Range('E23').Value = '  +16.96%  '
$ws.Range('E24').Value = '  +3.53%  '
$ws.Range('E26').Value = '  +3.15%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.90'
$ws.Range('E27').Value = '  +9.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.751.01'
$ws.Range('E28').Value = '  +5.49%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0948'
$ws.Range('E30').Value = '  +5.66%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '518.45'
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('E32').Value = '  +3.60%  '
$ws.Range('E33').Value = '  +6.49%  '
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.66'
$ws.Range('E36').Value = '  +2.33%  '
$ws.Range('E37').Value = '  +2.78%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.07'
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('E39').Value = '  +5.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.111'
$ws.Range('E40').Value = '  -6.50%  '
$ws.Range('E41').Value = '  +6.14%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('E43').Value = '  +4.95%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').Value = '  +9.12%  '
$ws.Range('E45').Value = '  +3.18%  '
$ws.Range('E46').Value = '  +1.66%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '150.07'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.541'
$ws.Range('E49').Value = '  +4.77%  '
$ws.Range('E50').Value = '  +7.22%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0261'
$ws.Range('E51').Value = '  +3.61%  '
